$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Value = "t1"
$ws.Range("Z1").Font.Size = 8
$ws.Range("Z2").Value = "t2"
$ws.Range("Z2").Font.Size = 8
Write-Host "done"
